# Apply "Atualização de bases das ligas, do dia: 11-04-2024 às 00:31"
# Greece Super League 1 sheet updates:
#  - Rows 194/195 swap their entire data (except column A, the row index)
#  - Rows 200/201 swap their entire data (except column A, the row index)
#  - Rows 206, 207, 209, 210, 211 get individual odds-column value corrections

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($Sheet, $Row, $Values) {
    foreach ($col in $Values.Keys) {
        $Sheet.Range("$col$Row").Value = $Values[$col]
    }
}

# ---- Row 194 becomes the old Row 195 content (column A stays 192) ----
Set-RowValues $ws 194 @{
    'B'  = 7920470
    'F'  = 'AEK Athens'
    'G'  = 'Olympiakos'
    'H'  = 1
    'I'  = 0
    'K'  = 1.909
    'L'  = 3.4
    'M'  = 4.2
    'N'  = 2.2
    'O'  = 3.2
    'P'  = 3.5
    'Q'  = -0.25
    'R'  = 1.85
    'S'  = 2
    'T'  = 2.5
    'W'  = 1.2
    'Z'  = 0.8500000000000001
    'AB' = -1
    'AC' = 0.825
}

# ---- Row 195 becomes the old Row 194 content (column A stays 193) ----
Set-RowValues $ws 195 @{
    'B'  = 7920471
    'F'  = 'Aris Salonika'
    'G'  = 'Lamia'
    'H'  = 3
    'I'  = 1
    'K'  = 1.571
    'L'  = 4
    'M'  = 6
    'N'  = 1.444
    'O'  = 4.5
    'P'  = 8.5
    'Q'  = -1.25
    'R'  = 1.925
    'S'  = 1.925
    'T'  = 2.75
    'W'  = 0.444
    'Z'  = 0.925
    'AB' = 1.025
    'AC' = -1
}

# ---- Row 200 becomes the old Row 201 content (column A stays 198) ----
Set-RowValues $ws 200 @{
    'B'  = 7920453
    'F'  = 'Panetolikos'
    'G'  = 'Volos NFC'
    'H'  = 0
    'I'  = 1
    'K'  = 2.3
    'L'  = 3
    'M'  = 3.4
    'O'  = 3.1
    'P'  = 3.8
    'R'  = 1.8
    'S'  = 2.05
    'T'  = 2.25
    'Y'  = 2.8
    'AA' = 1.05
    'AB' = -1
    'AC' = 0.825
}

# ---- Row 201 becomes the old Row 200 content (column A stays 199) ----
Set-RowValues $ws 201 @{
    'B'  = 7920450
    'F'  = 'Asteras Tripolis'
    'G'  = 'Kifisias FC'
    'H'  = 1
    'I'  = 2
    'K'  = 2.05
    'L'  = 3.3
    'M'  = 3.6
    'O'  = 3.5
    'P'  = 3.4
    'R'  = 1.85
    'S'  = 2
    'T'  = 2.75
    'Y'  = 2.4
    'AA' = 1
    'AB' = 0.5125
    'AC' = -0.5
}

# ---- Row 206: single odds correction ----
Set-RowValues $ws 206 @{
    'N' = 2.25
}

# ---- Row 207: odds corrections ----
Set-RowValues $ws 207 @{
    'N' = 2.05
    'P' = 3.8
    'U' = 1.875
    'V' = 1.975
}

# ---- Row 209: odds corrections ----
Set-RowValues $ws 209 @{
    'N' = 2.375
    'O' = 3
    'P' = 3.2
    'R' = 2.05
    'S' = 1.8
}

# ---- Row 210: odds corrections ----
Set-RowValues $ws 210 @{
    'U' = 2.025
    'V' = 1.825
}

# ---- Row 211: odds corrections ----
Set-RowValues $ws 211 @{
    'N' = 1.25
    'O' = 5.75
    'P' = 15
    'R' = 1.95
    'S' = 1.9
    'U' = 2.025
    'V' = 1.825
}
